$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (J1, K1, L1) - match the style used by the rest of row 1 (header)
$ws.Range("A1").Copy()
$ws.Range("J1:L1").PasteSpecial(-4122)
$ws.Range("J1").Value = "Duration"
$ws.Range("K1").Value = "ID"
$ws.Range("L1").Value = "Channel"

# Row 2 - video "不僅僅是編程｜如何成為成功的工程師" @ 00:06:39
$ws.Range("G2").Value = "職場"
$ws.Range("I2").Value = $true
$ws.Range("J2").Value = "7:51"
$ws.Range("K2").Value = "xS5Lv7-bMYI"
$ws.Range("L2").Value = "HackBear 泰瑞"

# Row 3 - video "不僅僅是編程｜如何成為成功的工程師" @ 00:07:17
$ws.Range("G3").Value = "職場"
$ws.Range("I3").Value = $true
$ws.Range("J3").Value = "7:51"
$ws.Range("K3").Value = "xS5Lv7-bMYI"
$ws.Range("L3").Value = "HackBear 泰瑞"

# Row 4 - video "不僅僅是編程｜如何成為成功的工程師" @ 00:07:19
$ws.Range("G4").Value = "職場"
$ws.Range("I4").Value = $true
$ws.Range("J4").Value = "7:51"
$ws.Range("K4").Value = "xS5Lv7-bMYI"
$ws.Range("L4").Value = "HackBear 泰瑞"

# Row 5 - video "不僅僅是編程｜如何成為成功的工程師" @ 00:07:22
$ws.Range("G5").Value = "職場"
$ws.Range("I5").Value = $true
$ws.Range("J5").Value = "7:51"
$ws.Range("K5").Value = "xS5Lv7-bMYI"
$ws.Range("L5").Value = "HackBear 泰瑞"

# Row 6 - video "好和弦教你做 8-bit 音樂！懷舊電玩風～" @ 00:02:44
$ws.Range("G6").Value = "編曲"
$ws.Range("I6").Value = $true
$ws.Range("J6").Value = "11:29"
$ws.Range("K6").Value = "8TbGCGDEgFk"
$ws.Range("L6").Value = "NiceChord (好和弦)"

# Row 7 - video "好和弦教你做 8-bit 音樂！懷舊電玩風～" @ 00:05:07
$ws.Range("G7").Value = "編曲"
$ws.Range("I7").Value = $true
$ws.Range("J7").Value = "11:29"
$ws.Range("K7").Value = "8TbGCGDEgFk"
$ws.Range("L7").Value = "NiceChord (好和弦)"

# Row 8 - video "超簡單又厲害的編曲軟體 - PixiTracker！" @ 00:10:02
$ws.Range("G8").Value = "編曲"
$ws.Range("I8").Value = $true
$ws.Range("J8").Value = "11:08"
$ws.Range("K8").Value = "hciSF-wGlyc"
$ws.Range("L8").Value = "NiceChord (好和弦)"

# Row 9 - video "超簡單又厲害的編曲軟體 - PixiTracker！" @ 00:10:04
$ws.Range("G9").Value = "編曲"
$ws.Range("I9").Value = $true
$ws.Range("J9").Value = "11:08"
$ws.Range("K9").Value = "hciSF-wGlyc"
$ws.Range("L9").Value = "NiceChord (好和弦)"
